$d = $word.ActiveDocument

# --- Update header date (single unique occurrence; whole-document find is safe) ---
$d.Content.Find.Execute("2023-08-17 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-08-18 Friday", 2) | Out-Null

# --- Update table cell values ---
# NOTE: Range.Find.Execute in this runtime operates over the whole document
# story regardless of the Range it is called on, so per-cell Find cannot be
# safely scoped (it would clobber the first document-order match, which can be
# in a different cell once earlier replacements change the text). Instead we
# assign Range.Text directly on each target cell, which replaces only that
# cell's content (formatting/end-of-cell markers are preserved) and is free
# of any cross-cell ordering hazards.
$tbl = $d.Tables.Item(1)

$tbl.Cell(1, 1).Range.Text = "37÷8=4, 5"
$tbl.Cell(1, 2).Range.Text = "38÷2=19, 0"
$tbl.Cell(1, 3).Range.Text = "79÷2=39, 1"
$tbl.Cell(1, 4).Range.Text = "29÷2=14, 1"
$tbl.Cell(1, 5).Range.Text = "71÷3=23, 2"
$tbl.Cell(5, 1).Range.Text = "83÷8=10, 3"
$tbl.Cell(5, 2).Range.Text = "60÷8=7, 4"
$tbl.Cell(5, 3).Range.Text = "77÷4=19, 1"
$tbl.Cell(5, 4).Range.Text = "39÷5=7, 4"
$tbl.Cell(5, 5).Range.Text = "42÷6=7, 0"
$tbl.Cell(9, 1).Range.Text = "45÷9=5, 0"
$tbl.Cell(9, 2).Range.Text = "59÷4=14, 3"
$tbl.Cell(9, 3).Range.Text = "36÷9=4, 0"
$tbl.Cell(9, 4).Range.Text = "38÷2=19, 0"
$tbl.Cell(9, 5).Range.Text = "37÷3=12, 1"
$tbl.Cell(13, 1).Range.Text = "47÷9=5, 2"
$tbl.Cell(13, 2).Range.Text = "13÷5=2, 3"
$tbl.Cell(13, 3).Range.Text = "38÷7=5, 3"
$tbl.Cell(13, 4).Range.Text = "59÷4=14, 3"
$tbl.Cell(13, 5).Range.Text = "19÷9=2, 1"
$tbl.Cell(17, 1).Range.Text = "35÷8=4, 3"
$tbl.Cell(17, 2).Range.Text = "21÷4=5, 1"
$tbl.Cell(17, 3).Range.Text = "39÷4=9, 3"
$tbl.Cell(17, 4).Range.Text = "10÷8=1, 2"
$tbl.Cell(17, 5).Range.Text = "62÷2=31, 0"

Write-Output "done"
